$wb = $excel.ActiveWorkbook
$objecten = $wb.Worksheets.Item("Objecten")
$parameters = $wb.Worksheets.Item("Parameters")

# Add the new "ClassificatieCode" column (B) on the Objecten sheet, next to
# the existing "Onderdeel" column (A).
$objecten.Range("B1").Value = "ClassificatieCode"
$objecten.Range("B2").Value = "BM"
$objecten.Range("B3").Value = "SI"
$objecten.Range("B4").Value = "RE"
$objecten.Range("B5").Value = "RF"
$objecten.Range("B6").Value = "SE"
$objecten.Range("B7").Value = "IB"
$objecten.Range("B8").Value = "VB"
$objecten.Range("B9").Value = "TI"

# Set the new column's width to match the source workbook (closest value
# the engine's character-width rounding can reach to the target 17.63).
$objecten.Columns.Item(2).ColumnWidth = 16.8

# Restore the default selection on the Parameters sheet.
$parameters.Range("A2").Select() | Out-Null

# Update the selection on the Objecten sheet.
$objecten.Range("B6").Select() | Out-Null

# Make "Objecten" the active/selected sheet instead of "Parameters" (must be
# last so it is the one left active/selected in the saved workbook).
$objecten.Activate() | Out-Null
